$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (this shifts existing rows 19-96 down to 20-97)
$ws.Rows(19).Insert()

# Populate the newly inserted row 19 with the weekly data entry
# (same as the prior top row, but with an updated "Fecha" / date value)
$ws.Cells.Item(19, 1).Value2 = 9
$ws.Cells.Item(19, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value2 = "Metropolitana"
$ws.Cells.Item(19, 4).Value2 = 44676
$ws.Cells.Item(19, 5).Value2 = 13
$ws.Cells.Item(19, 6).Value2 = 100114007
$ws.Cells.Item(19, 7).Value2 = "Jengibre"
$ws.Cells.Item(19, 8).Value2 = "Sin especificar"
$ws.Cells.Item(19, 9).Value2 = "Primera"
$ws.Cells.Item(19, 10).Value2 = 790
$ws.Cells.Item(19, 11).Value2 = 11000
$ws.Cells.Item(19, 12).Value2 = 12000
$ws.Cells.Item(19, 13).Value2 = 11494
$ws.Cells.Item(19, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(19, 15).Value2 = "Perú"
$ws.Cells.Item(19, 16).Value2 = 884
$ws.Cells.Item(19, 17).Value2 = 13
$ws.Cells.Item(19, 18).Value2 = "Hortaliza"

# Match the date cell's number format style used by the rest of the "Fecha" column
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
